$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 2966.5715
$ws.Range("I113").Value = 1824
$ws.Range("K113").Value = 1824
$ws.Range("M113").Value = 1430
$ws.Range("H116").Value = 3235
$ws.Range("I116").Value = 4735
$ws.Range("J116").Value = 2672.5
$ws.Range("K116").Value = 4735
$ws.Range("L116").Value = 2672.5
$ws.Range("M116").Value = -1293
$ws.Range("N116").Value = -9556.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 46293.727
$ws.Range("I2").Value = 59597.883
$ws.Range("J2").Value = 1059.6
$ws.Range("K2").Value = 59597.883
$ws.Range("L2").Value = 1059.6
$ws.Range("M2").Value = -59484.883
$ws.Range("N2").Value = -1285.6
$ws.Range("H32").Value = 3345.51
$ws.Range("I32").Value = 2335.4175
$ws.Range("J32").Value = 13558.667
$ws.Range("K32").Value = 2335.4175
$ws.Range("L32").Value = 13558.667
$ws.Range("M32").Value = -2048.4175
$ws.Range("N32").Value = -14132.667
$ws.Range("H116").Value = 46293.727
$ws.Range("I116").Value = 59597.883
$ws.Range("J116").Value = 1059.6
$ws.Range("K116").Value = 59597.883
$ws.Range("L116").Value = 1059.6
$ws.Range("M116").Value = -57303.883
$ws.Range("N116").Value = -5647.6
$ws.Range("H132").Value = 8103.161
$ws.Range("I132").Value = 5244.5386
$ws.Range("K132").Value = 15733.6158
$ws.Range("M132").Value = -13203.6158
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 46293.727
$ws.Range("I3").Value = 59597.883
$ws.Range("J3").Value = 1059.6
$ws.Range("K3").Value = 59597.883
$ws.Range("L3").Value = 1059.6
$ws.Range("M3").Value = -59483.883
$ws.Range("N3").Value = -1287.6
$ws.Range("H86").Value = 1829.8158
$ws.Range("I86").Value = 1602.75
$ws.Range("J86").Value = 2082.111
$ws.Range("K86").Value = 1602.75
$ws.Range("L86").Value = 2082.111
$ws.Range("M86").Value = -479.75
$ws.Range("N86").Value = -4328.111
$ws.Range("H89").Value = 1829.8158
$ws.Range("I89").Value = 1602.75
$ws.Range("J89").Value = 2082.111
$ws.Range("K89").Value = 8013.75
$ws.Range("L89").Value = 10410.555
$ws.Range("M89").Value = -2397.75
$ws.Range("N89").Value = -21642.555
$ws.Range("H105").Value = 1537.4193
$ws.Range("I105").Value = 1218.9524
$ws.Range("K105").Value = 1218.9524
$ws.Range("M105").Value = 528.0476000000001
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 948.5345
$ws.Range("I58").Value = 621.87177
$ws.Range("K58").Value = 621.87177
$ws.Range("M58").Value = -418.87177
$ws.Range("H105").Value = 1713.129
$ws.Range("I105").Value = 1102.1364
$ws.Range("J105").Value = 3206.6667
$ws.Range("K105").Value = 1102.1364
$ws.Range("L105").Value = 3206.6667
$ws.Range("M105").Value = 644.8635999999999
$ws.Range("N105").Value = -6700.6667
$ws.Range("H134").Value = 2226.5454
$ws.Range("I134").Value = 2209.7036
$ws.Range("J134").Value = 2302.3333
$ws.Range("K134").Value = 6629.110799999999
$ws.Range("L134").Value = 6906.999899999999
$ws.Range("M134").Value = -4094.110799999999
$ws.Range("N134").Value = -11976.9999
$ws.Range("H136").Value = 948.5345
$ws.Range("I136").Value = 621.87177
$ws.Range("K136").Value = 1865.61531
$ws.Range("M136").Value = 684.3846900000001
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 38.64516
$ws.Range("I12").Value = 2.5
$ws.Range("J12").Value = 47.32
$ws.Range("K12").Value = 7.5
$ws.Range("L12").Value = 141.96
$ws.Range("M12").Value = 165.5
$ws.Range("N12").Value = -487.96
$ws.Range("H38").Value = 608.6667
$ws.Range("I38").Value = 31.6
$ws.Range("J38").Value = 830.61536
$ws.Range("K38").Value = 94.80000000000001
$ws.Range("L38").Value = 2491.84608
$ws.Range("M38").Value = 252.2
$ws.Range("N38").Value = -3185.84608
$ws.Range("H62").Value = 3125
$ws.Range("J62").Value = 3125
$ws.Range("L62").Value = 9375
$ws.Range("N62").Value = -10747
$ws.Range("H65").Value = 3125
$ws.Range("J65").Value = 3125
$ws.Range("L65").Value = 28125
$ws.Range("N65").Value = -34989
$ws.Range("H98").Value = 385.45
$ws.Range("I98").Value = 294.9375
$ws.Range("J98").Value = 747.5
$ws.Range("K98").Value = 884.8125
$ws.Range("L98").Value = 2242.5
$ws.Range("M98").Value = 613.1875
$ws.Range("N98").Value = -5238.5
$ws.Range("H139").Value = 335098.7
$ws.Range("I139").Value = 459756.6
$ws.Range("K139").Value = 1379269.8
$ws.Range("M139").Value = -1374129.8
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 215.8
$ws.Range("I107").Value = 188.5
$ws.Range("J107").Value = 325
$ws.Range("K107").Value = 188.5
$ws.Range("L107").Value = 325
$ws.Range("M107").Value = 1731.5
$ws.Range("N107").Value = -4165
$ws.Range("H126").Value = 1454.4546
$ws.Range("I126").Value = 1120
$ws.Range("J126").Value = 1733.1666
$ws.Range("K126").Value = 3360
$ws.Range("L126").Value = 5199.4998
$ws.Range("M126").Value = -890
$ws.Range("N126").Value = -10139.4998
$ws.Range("H132").Value = 7363.7827
$ws.Range("I132").Value = 9134.625
$ws.Range("J132").Value = 3316.1428
$ws.Range("K132").Value = 27403.875
$ws.Range("L132").Value = 9948.428400000001
$ws.Range("M132").Value = -24873.875
$ws.Range("N132").Value = -15008.4284
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1640.909
$ws.Range("I61").Value = 1360
$ws.Range("J61").Value = 1823.5
$ws.Range("K61").Value = 1360
$ws.Range("L61").Value = 1823.5
$ws.Range("M61").Value = -1158
$ws.Range("N61").Value = -2227.5
$ws.Range("H107").Value = 2374.75
$ws.Range("I107").Value = 2374.75
$ws.Range("K107").Value = 2374.75
$ws.Range("M107").Value = -454.75
$ws.Range("H113").Value = 1640.909
$ws.Range("I113").Value = 1360
$ws.Range("J113").Value = 1823.5
$ws.Range("K113").Value = 1360
$ws.Range("L113").Value = 1823.5
$ws.Range("M113").Value = 810
$ws.Range("N113").Value = -6163.5
$ws.Range("H132").Value = 3925
$ws.Range("I132").Value = 4405.6055
$ws.Range("J132").Value = 3094.8635
$ws.Range("K132").Value = 13216.8165
$ws.Range("L132").Value = 9284.5905
$ws.Range("M132").Value = -10686.8165
$ws.Range("N132").Value = -14344.5905
$ws.Range("H136").Value = 3201.8086
$ws.Range("I136").Value = 1037.8889
$ws.Range("J136").Value = 51890
$ws.Range("K136").Value = 3113.6667
$ws.Range("L136").Value = 155670
$ws.Range("M136").Value = -563.6666999999998
$ws.Range("N136").Value = -160770
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 664.36365
$ws.Range("I100").Value = 758.8333
$ws.Range("J100").Value = 551
$ws.Range("K100").Value = 1517.6666
$ws.Range("L100").Value = 1102
$ws.Range("M100").Value = -976.6666
$ws.Range("N100").Value = -2184
$ws.Range("H132").Value = 7578307.5
$ws.Range("I132").Value = 11907932
$ws.Range("J132").Value = 1463.75
$ws.Range("K132").Value = 35723796
$ws.Range("L132").Value = 4391.25
$ws.Range("M132").Value = -35721266
$ws.Range("N132").Value = -9451.25
